# Applies the "hour log" update described in the commit:
#   - Added additional search filters (+ row for pagination attempt)
#   - Added a "Done" wrap-up row
#   - Updated SRP description with a reference link
#   - Added a hyperlink on the pagination row pointing to a YouTube video
#   - Scrolled the sheet back to the top (A2)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the SRP description (row 14, column D) with the extra sentence ---
$ws.Range("D14").Value = "I didn't really looked into ""SRP"" (Signle responsibilty protocol) when i saw it mentioned in the assessment, but when i went through the checklist i saw it come by, i decided to look into it. This turned out to be a very time consuming choice because it took alot of time to refactor the code. I found this website (https://medium.com/@Omojunior11/single-responsibility-principle-srp-example-using-php-337e33d739e) very usefull for understanding the basics. I added a new services method with 2 services. These are then used in the controller. This way the, for example, the methods that validate input become reusable and more maintainable. I added a __construct in the FacilityController to initiate the DB connection and creating instances of the services, which then get the `$db variable injected. I have seen code being structured like this but i never really understood why it was now."

# --- Row 16: new log entry about pagination / search filters ---
$ws.Range("A16").Value = "I decided to try to add pagnation for the read operations and added new search filters to the search operation"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = Get-Date -Year 2025 -Month 4 -Day 29
$ws.Range("D16").Value = "After reading through the assessment once more, i decided to try to add another optional feature, the pagination for read operations. I have never done this before so i looked up some tutorials on the internet about it. I used this blog (https://www.merge.dev/blog/rest-api-pagination) and AI to understand how pagination works and how to implement it. After a while i could not figure out how to make it work so i went back to the way it was. I need to know pagination better before using it in a pretty complex API like this one. I did add some additional search filters to the search operation though!"
$ws.Range("A16").RowHeight = 39.75

# --- Row 17: wrap-up entry ---
$ws.Range("A17").Value = "Done"
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = Get-Date -Year 2025 -Month 4 -Day 29
$ws.Range("D17").Value = "I think im ready to finish this assessment. I spent 3 days working on this assessment as much and hard as i could. I believe i included all requirements and did everything that needed to be done for the assessment. I had an absolute blast making this thing, because i enjoyed learning / making something that i know will be important to do with my eyes closed in the future. I have learned a lot of new things, which im very excited about. Thank you guys."

# --- Hyperlink on D16 pointing to the referenced video ---
$ws.Hyperlinks.Add($ws.Range("D16"), "https://youtu.be/KrHkOzJxhss", "", "", "https://youtu.be/KrHkOzJxhss")

# --- Row 15 picks up the taller auto-fit height used by the neighbouring rows ---
$ws.Range("A15").RowHeight = 39.75

# --- Scroll position: back to the top-left of the sheet ---
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Application.ActiveWindow.ScrollColumn = 1

$wb.Save()
